$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.715.89"
$ws.Range("E2").Value = "  +1.50%  "
$ws.Range("D3").Value = "1.565.40"
$ws.Range("E3").Value = "  +0.05%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.31%  "
$ws.Range("D5").Value = "'210.08"
$ws.Range("E5").Value = "  -0.23%  "
$ws.Range("E6").Value = "  +2.12%  "
$ws.Range("E7").Value = "  -0.44%  "
$ws.Range("D8").Value = "'25.10"
$ws.Range("E8").Value = "  +5.67%  "
$ws.Range("D9").Value = "'0.245"
$ws.Range("E9").Value = "  +0.71%  "
$ws.Range("D10").Value = "'0.0587"
$ws.Range("E10").Value = "  +0.29%  "
$ws.Range("D11").Value = "'0.0896"
$ws.Range("E11").Value = "  +0.13%  "
$ws.Range("D12").Value = "1.791.32"
$ws.Range("E12").Value = "  +0.15%  "
$ws.Range("D13").Value = "1.563.98"
$ws.Range("E13").Value = "  -0.11%  "
$ws.Range("D14").Value = "28.729.64"
$ws.Range("E14").Value = "  +1.57%  "
$ws.Range("E15").Value = "  +0.89%  "
$ws.Range("E16").Value = "  -0.24%  "
$ws.Range("E17").Value = "  +0.80%  "
$ws.Range("D18").Value = "'230.33"
$ws.Range("E18").Value = "  +1.18%  "
$ws.Range("D19").Value = "'7.35"
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("E20").Value = "  +0.73%  "
$ws.Range("E21").Value = "  -0.49%  "
$ws.Range("D22").Value = "'3.92"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").Value = "'9.05"
$ws.Range("E23").Value = "  +1.66%  "
$ws.Range("E24").Value = "  +2.10%  "
$ws.Range("D25").Value = "'151.37"
$ws.Range("E25").Value = "  +0.76%  "
$ws.Range("D26").Value = "'14.80"
$ws.Range("E26").Value = "  -0.40%  "
$ws.Range("E27").Value = "  +1.78%  "
$ws.Range("D28").Value = "'0.999"
$ws.Range("E28").Value = "  -0.39%  "
$ws.Range("E29").Value = "  -1.34%  "
$ws.Range("E30").Value = "  -3.56%  "
$ws.Range("E31").Value = "  -2.01%  "
$ws.Range("D32").Value = "'3.18"
$ws.Range("E32").Value = "  +0.53%  "
$ws.Range("D33").Value = "1.398.11"
$ws.Range("D34").Value = "'2.99"
$ws.Range("E34").Value = "  -2.69%  "
$ws.Range("E35").Value = "  -3.52%  "
$ws.Range("D36").Value = "'1.47"
$ws.Range("E36").Value = "  -1.46%  "
$ws.Range("D37").Value = "'2.69"
$ws.Range("E37").Value = "  +2.06%  "
$ws.Range("E38").Value = "  -2.35%  "
$ws.Range("E39").Value = "  -0.43%  "
$ws.Range("D40").Value = "'1.95"
$ws.Range("E40").Value = "  +1.17%  "
$ws.Range("D41").Value = "'0.518"
$ws.Range("E41").Value = "  -0.27%  "
$ws.Range("D42").Value = "'0.997"
$ws.Range("E42").Value = "  -0.44%  "
$ws.Range("D43").Value = "'0.772"
$ws.Range("E43").Value = "  -1.08%  "
$ws.Range("E44").Value = "  -3.64%  "
$ws.Range("D45").Value = "'63.95"
$ws.Range("E45").Value = "  +2.98%  "
$ws.Range("E46").Value = "  -1.45%  "
$ws.Range("D47").Value = "1.703.34"
$ws.Range("E47").Value = "  +0.09%  "
$ws.Range("D48").Value = "'0.870"
$ws.Range("E48").Value = "  -5.39%  "
$ws.Range("E49").Value = "  -0.29%  "
$ws.Range("D50").Value = "'42.71"
$ws.Range("E50").Value = "  +5.22%  "
$ws.Range("D51").Value = "'0.0512"
$ws.Range("E51").Value = "  -0.38%  "
